$d = $word.ActiveDocument

# --- Edit 1: add a new sentence (as a separate run) after the "texture class dll." note ---
$targetPara1 = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    if ($para.Range.Text -like "*texture class dll.*") {
        $targetPara1 = $para
        break
    }
}
if ($targetPara1 -ne $null) {
    $r1 = $targetPara1.Range
    $r1.MoveEnd(1, -1)
    $r1.Collapse(0)
    $r1.InsertAfter("But, I made a modification to the vba code to use a constant texture class for a solute diffusion parameter. This forces it to work in 64 bit and the model is not very sensitive to this value.")
}

# --- Edit 2: split "Read Plant filesV6.xlsm" (the standalone file-listing paragraph) into
#     three runs: "Read Plant filesV" / "7_mulch" / ".xlsm" ---
$targetPara2 = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    if ($para.Range.Text.TrimEnd() -eq "Read Plant filesV6.xlsm") {
        $targetPara2 = $para
        break
    }
}
if ($targetPara2 -ne $null) {
    $full2 = $targetPara2.Range
    $full2.MoveEnd(1, -1)
    $start2 = $full2.Start
    # "Read Plant filesV" is 17 characters long; the "6" sits right after it.
    $rVersion = $d.Range($start2 + 17, $start2 + 18)

    # Use tracked changes so the replacement lands as its own run instead of being
    # folded back into the neighboring plain-text run, then accept just those
    # revisions individually (not Document.AcceptAllRevisions, which also touches
    # unrelated rsid bookkeeping elsewhere in the document) to flatten the markup
    # away while keeping the new run boundaries.
    $wasTracking = $d.TrackRevisions
    $d.TrackRevisions = $true
    $rVersion.Text = "7_mulch"
    $d.TrackRevisions = $wasTracking
    for ($j = $d.Revisions.Count; $j -ge 1; $j--) {
        $d.Revisions.Item($j).Accept()
    }
}
